$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("T1").Value = 0.1892678598946218
$ws.Range("T2").Value = 0.032302043749510331
$ws.Range("T3").Value = 0.026238275402749485
$ws.Range("T4").Value = 0.16905002841257127
$ws.Range("T5").Value = 0.19504265042326546
$ws.Range("T6").Value = 0.15611659276390244
$ws.Range("T7").Value = 0.10617299893863949
$ws.Range("T8").Value = 0.061919466480956781
$ws.Range("T9").Value = 0.061919466480956781
$ws.Range("T10").Value = 0.0019706174528261178
